$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.105.44"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.838.88"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.71%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6286"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.10%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07584"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.38%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2936"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.61"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.77%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07757"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.841.14"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.77%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.975"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6666"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "83.25"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001000"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +15.51%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.075"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.71%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.125.48"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "227.24"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.19%  "
$ws.Range("E22").Value = "  +1.62%  "
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "159.99"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.22%  "
$ws.Range("E25").Value = "  +1.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.515"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.60%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.97"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.498"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.108"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.45%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.019"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.20%  "
$ws.Range("E31").Value = "  -0.50%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.05267"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.848"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.54%  "
$ws.Range("E34").Value = "  +0.23%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.138"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.683"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.29%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.245.17"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.89%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.765"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01788"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.362"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9029"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.63%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "102.10"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.00000000125"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.987.94"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.70%  "
$ws.Range("E46").Value = "  +0.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5129"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4047"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.929"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.61%  "
$ws.Range("E50").Value = "  -0.54%  "
$ws.Range("E51").Value = "  +0.45%  "
